$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header cell G1 ("elementIndex"), formatted like the other header cells
$ws.Range("F1").Copy()
$ws.Range("G1").PasteSpecial(-4122)
$ws.Range("G1").Value = "elementIndex"

# New data column G2:G9 holding the numeric element index (0-based), formatted
# like the rest of the data rows
for ($i = 0; $i -le 7; $i++) {
    $row = 2 + $i
    $srcCell = $ws.Cells.Item($row, 6)
    $dstCell = $ws.Cells.Item($row, 7)
    $srcCell.Copy()
    $dstCell.PasteSpecial(-4122)
    $dstCell.Value = $i
}

$excel.CutCopyMode = 0
